$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the (second) blank paragraph that sits directly above the
#    "Date: " paragraph into the "Date: " paragraph, i.e. delete that
#    blank paragraph's end-of-paragraph mark.
# ---------------------------------------------------------------------
$dateRange = $d.Content
$dateRange.Find.Execute("Date: ", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$datePara = $dateRange.Paragraphs.Item(1)
$mark = $d.Range($datePara.Range.Start - 1, $datePara.Range.Start)
$mark.Delete()

# ---------------------------------------------------------------------
# 2. Remove the word "Apple" that follows "Device: ", leaving just the
#    label text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Device: Apple", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Device: ", 2)

# ---------------------------------------------------------------------
# 3. Add a new paragraph right after "User agent: " containing the text
#    "Client App used:".
# ---------------------------------------------------------------------
$uaRange = $d.Content
$uaRange.Find.Execute("User agent: ", $false, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)
$uaPara = $uaRange.Paragraphs.Item(1)
$uaParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $uaPara.Range.Start) {
        $uaParaIndex = $i
        break
    }
}
$uaPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($uaParaIndex + 1)
$insertRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertRange.InsertAfter("Client App u")
$newPara2 = $d.Paragraphs.Item($uaParaIndex + 1)
$pos = $newPara2.Range.End - 1
$tailRange = $d.Range($pos, $pos)
$tailRange.InsertAfter("sed:")
